# Update the default model data for Imperial, shifting dates forward
# by 28 days (25.06.2024 -> 23.07.2024) and refreshing the Prediction
# column (C) with the newer model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dayShift = 28
$newDateLabel = "23.07.2024"

# New Prediction (column C) values, keyed by the Interval (column B) value.
$predictions = @{
    1  = 0
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 0.014
    30 = 0.032
    31 = 0.056
    32 = 0.087
    33 = 0.124
    34 = 0.166
    35 = 0.21
    36 = 0.255
    37 = 0.298
    38 = 0.342
    39 = 0.383
    40 = 0.449
    41 = 0.52
    42 = 0.584
    43 = 0.616
    44 = 0.641
    45 = 0.68
    46 = 0.697
    47 = 0.706
    48 = 0.717
    49 = 0.732
    50 = 0.733
    51 = 0.733
    52 = 0.733
    53 = 0.727
    54 = 0.727
    55 = 0.724
    56 = 0.715
    57 = 0.71
    58 = 0.703
    59 = 0.683
    60 = 0.672
    61 = 0.672
    62 = 0.662
    63 = 0.622
    64 = 0.587
    65 = 0.554
    66 = 0.534
    67 = 0.501
    68 = 0.458
    69 = 0.416
    70 = 0.37
    71 = 0.324
    72 = 0.307
    73 = 0.283
    74 = 0.257
    75 = 0.212
    76 = 0.163
    77 = 0.144
    78 = 0.117
    79 = 0.098
    80 = 0.086
    81 = 0.067
    82 = 0.052
    83 = 0.041
    84 = 0.033
    85 = 0.024
    86 = 0
    87 = 0
    88 = 0
    89 = 0
    90 = 0
    91 = 0
    92 = 0
    93 = 0
    94 = 0
    95 = 0
    96 = 0
}

for ($r = 2; $r -le 96; $r++) {
    $interval = $ws.Cells.Item($r, 2).Value2()

    # Column A: shift the date/time serial forward by 28 days, keeping
    # the same time-of-day fraction and number format.
    $oldSerial = $ws.Cells.Item($r, 1).Value2()
    $ws.Cells.Item($r, 1).Value = $oldSerial + $dayShift

    # Column C: refresh with the new prediction value.
    $ws.Cells.Item($r, 3).Value = $predictions[[int]$interval]

    # Column D: rebuild the lookup label with the new date.
    $ws.Cells.Item($r, 4).Value = "$newDateLabel$interval"
}
